$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in header template string: "tocketsInWork" -> "ticketsInWork"
$ws.Range("E1").Value = "{d.i18n.ticketsInWork}"
